$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 259
$ws.Range("I2").Value = 669
$ws.Range("J2").Value = 2817
$ws.Range("K2").Value = 17
$ws.Range("L2").Value = 760
$ws.Range("M2").Value = 34
$ws.Range("N2").Value = 470
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 11
$ws.Range("S2").Value = 288
$ws.Range("T2").Value = 481
$ws.Range("U2").Value = 31
$ws.Range("V2").Value = 4351
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 4294
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 67
$ws.Range("AA2").Value = 24
